$wb = $excel.ActiveWorkbook

# ---- classFields sheet (xl/worksheets/sheet3.xml) ----
# Field rows reordered/regrouped per class (e.g. LOG/SOURCE ordering,
# "id" moved after the other Customer fields, etc.)
$ws3 = $wb.Worksheets.Item("classFields")

$ws3.Cells.Item(4,1).Value = 'pl.piomin.payment.service.OrderManageService'
$ws3.Cells.Item(4,2).Value = 'SOURCE'
$ws3.Cells.Item(4,3).Value = 'private'
$ws3.Cells.Item(4,4).Value = 'java.lang.String'

$ws3.Cells.Item(5,1).Value = 'pl.piomin.payment.service.OrderManageService'
$ws3.Cells.Item(5,2).Value = 'LOG'
$ws3.Cells.Item(5,3).Value = 'private'
$ws3.Cells.Item(5,4).Value = 'org.slf4j.Logger'

$ws3.Cells.Item(6,1).Value = 'pl.piomin.payment.PaymentComponentTests'
$ws3.Cells.Item(6,2).Value = 'LOG'
$ws3.Cells.Item(6,3).Value = 'private'
$ws3.Cells.Item(6,4).Value = 'org.slf4j.Logger'

$ws3.Cells.Item(7,1).Value = 'pl.piomin.payment.PaymentComponentTests'
$ws3.Cells.Item(7,2).Value = 'kafka'
$ws3.Cells.Item(7,3).Value = 'private'
$ws3.Cells.Item(7,4).Value = 'org.springframework.kafka.test.EmbeddedKafkaBroker'

$ws3.Cells.Item(8,1).Value = 'pl.piomin.payment.PaymentComponentTests'
$ws3.Cells.Item(8,2).Value = 'template'
$ws3.Cells.Item(8,3).Value = 'private'
$ws3.Cells.Item(8,4).Value = 'org.springframework.kafka.core.KafkaTemplate'

$ws3.Cells.Item(9,1).Value = 'pl.piomin.payment.PaymentComponentTests'
$ws3.Cells.Item(9,2).Value = 'factory'
$ws3.Cells.Item(9,3).Value = 'private'
$ws3.Cells.Item(9,4).Value = 'org.springframework.kafka.core.ConsumerFactory'

$ws3.Cells.Item(10,1).Value = 'pl.piomin.payment.PaymentComponentTests'
$ws3.Cells.Item(10,2).Value = 'repository'
$ws3.Cells.Item(10,3).Value = ''
$ws3.Cells.Item(10,4).Value = 'pl.piomin.payment.repository.CustomerRepository'

$ws3.Cells.Item(11,1).Value = 'pl.piomin.payment.PaymentComponentTests'
$ws3.Cells.Item(11,2).Value = 'customer'
$ws3.Cells.Item(11,3).Value = ''
$ws3.Cells.Item(11,4).Value = 'pl.piomin.payment.domain.Customer'

$ws3.Cells.Item(12,1).Value = 'pl.piomin.payment.PaymentApp'
$ws3.Cells.Item(12,2).Value = 'LOG'
$ws3.Cells.Item(12,3).Value = 'private'
$ws3.Cells.Item(12,4).Value = 'org.slf4j.Logger'

$ws3.Cells.Item(13,1).Value = 'pl.piomin.payment.PaymentApp'
$ws3.Cells.Item(13,2).Value = 'orderManageService'
$ws3.Cells.Item(13,3).Value = ''
$ws3.Cells.Item(13,4).Value = 'pl.piomin.payment.service.OrderManageService'

$ws3.Cells.Item(14,1).Value = 'pl.piomin.payment.PaymentApp'
$ws3.Cells.Item(14,2).Value = 'repository'
$ws3.Cells.Item(14,3).Value = 'private'
$ws3.Cells.Item(14,4).Value = 'pl.piomin.payment.repository.CustomerRepository'

$ws3.Cells.Item(15,1).Value = 'pl.piomin.payment.domain.Customer'
$ws3.Cells.Item(15,2).Value = 'amountAvailable'
$ws3.Cells.Item(15,3).Value = 'private'
$ws3.Cells.Item(15,4).Value = 'int'

$ws3.Cells.Item(16,1).Value = 'pl.piomin.payment.domain.Customer'
$ws3.Cells.Item(16,2).Value = 'amountReserved'
$ws3.Cells.Item(16,3).Value = 'private'
$ws3.Cells.Item(16,4).Value = 'int'

$ws3.Cells.Item(17,1).Value = 'pl.piomin.payment.domain.Customer'
$ws3.Cells.Item(17,2).Value = 'name'
$ws3.Cells.Item(17,3).Value = 'private'
$ws3.Cells.Item(17,4).Value = 'java.lang.String'

$ws3.Cells.Item(18,1).Value = 'pl.piomin.payment.domain.Customer'
$ws3.Cells.Item(18,2).Value = 'id'
$ws3.Cells.Item(18,3).Value = 'private'
$ws3.Cells.Item(18,4).Value = 'java.lang.Long'


# ---- methodNumberOfLines sheet (xl/worksheets/sheet11.xml) ----
# Added constructor rows (common-package handling) + reordered existing rows.
$ws11 = $wb.Worksheets.Item("methodNumberOfLines")

# Column C ("Number of Lines") holds numeric-looking text; pre-format as Text
# so values are written as strings (matching the source data type) instead of
# being auto-converted to numbers by Excel.
$ws11.Range("C2:C25").NumberFormat = "@"

$ws11.Cells.Item(2,1).Value = 'pl.piomin.payment.service.OrderManageService'
$ws11.Cells.Item(2,2).Value = 'OrderManageService(pl.piomin.payment.repository.CustomerRepository, org.springframework.kafka.core.KafkaTemplate)'
$ws11.Cells.Item(2,3).Value = '4'

$ws11.Cells.Item(3,1).Value = 'pl.piomin.payment.service.OrderManageService'
$ws11.Cells.Item(3,2).Value = 'reserve(pl.piomin.base.domain.Order)'
$ws11.Cells.Item(3,3).Value = '16'

$ws11.Cells.Item(4,1).Value = 'pl.piomin.payment.service.OrderManageService'
$ws11.Cells.Item(4,2).Value = 'confirm(pl.piomin.base.domain.Order)'
$ws11.Cells.Item(4,3).Value = '13'

$ws11.Cells.Item(5,1).Value = 'pl.piomin.payment.PaymentAppTest'
$ws11.Cells.Item(5,2).Value = 'PaymentAppTest()'
$ws11.Cells.Item(5,3).Value = '1'

$ws11.Cells.Item(6,1).Value = 'pl.piomin.payment.PaymentAppTest'
$ws11.Cells.Item(6,2).Value = 'main(java.lang.String[])'
$ws11.Cells.Item(6,3).Value = '3'

$ws11.Cells.Item(7,1).Value = 'pl.piomin.payment.PaymentComponentTests'
$ws11.Cells.Item(7,2).Value = 'PaymentComponentTests()'
$ws11.Cells.Item(7,3).Value = '1'

$ws11.Cells.Item(8,1).Value = 'pl.piomin.payment.PaymentComponentTests'
$ws11.Cells.Item(8,2).Value = 'eventAccept()'
$ws11.Cells.Item(8,3).Value = '11'

$ws11.Cells.Item(9,1).Value = 'pl.piomin.payment.PaymentComponentTests'
$ws11.Cells.Item(9,2).Value = 'eventReject()'
$ws11.Cells.Item(9,3).Value = '10'

$ws11.Cells.Item(10,1).Value = 'pl.piomin.payment.PaymentComponentTests'
$ws11.Cells.Item(10,2).Value = 'eventConfirm()'
$ws11.Cells.Item(10,3).Value = '10'

$ws11.Cells.Item(11,1).Value = 'pl.piomin.payment.PaymentApp'
$ws11.Cells.Item(11,2).Value = 'PaymentApp()'
$ws11.Cells.Item(11,3).Value = '1'

$ws11.Cells.Item(12,1).Value = 'pl.piomin.payment.PaymentApp'
$ws11.Cells.Item(12,2).Value = 'main(java.lang.String[])'
$ws11.Cells.Item(12,3).Value = '3'

$ws11.Cells.Item(13,1).Value = 'pl.piomin.payment.PaymentApp'
$ws11.Cells.Item(13,2).Value = 'onEvent(pl.piomin.base.domain.Order)'
$ws11.Cells.Item(13,3).Value = '9'

$ws11.Cells.Item(14,1).Value = 'pl.piomin.payment.PaymentApp'
$ws11.Cells.Item(14,2).Value = 'generateData()'
$ws11.Cells.Item(14,3).Value = '9'

$ws11.Cells.Item(15,1).Value = 'pl.piomin.payment.domain.Customer'
$ws11.Cells.Item(15,2).Value = 'getId()'
$ws11.Cells.Item(15,3).Value = '3'

$ws11.Cells.Item(16,1).Value = 'pl.piomin.payment.domain.Customer'
$ws11.Cells.Item(16,2).Value = 'setId(java.lang.Long)'
$ws11.Cells.Item(16,3).Value = '3'

$ws11.Cells.Item(17,1).Value = 'pl.piomin.payment.domain.Customer'
$ws11.Cells.Item(17,2).Value = 'getName()'
$ws11.Cells.Item(17,3).Value = '3'

$ws11.Cells.Item(18,1).Value = 'pl.piomin.payment.domain.Customer'
$ws11.Cells.Item(18,2).Value = 'setName(java.lang.String)'
$ws11.Cells.Item(18,3).Value = '3'

$ws11.Cells.Item(19,1).Value = 'pl.piomin.payment.domain.Customer'
$ws11.Cells.Item(19,2).Value = 'getAmountAvailable()'
$ws11.Cells.Item(19,3).Value = '3'

$ws11.Cells.Item(20,1).Value = 'pl.piomin.payment.domain.Customer'
$ws11.Cells.Item(20,2).Value = 'setAmountAvailable(int)'
$ws11.Cells.Item(20,3).Value = '3'

$ws11.Cells.Item(21,1).Value = 'pl.piomin.payment.domain.Customer'
$ws11.Cells.Item(21,2).Value = 'getAmountReserved()'
$ws11.Cells.Item(21,3).Value = '3'

$ws11.Cells.Item(22,1).Value = 'pl.piomin.payment.domain.Customer'
$ws11.Cells.Item(22,2).Value = 'setAmountReserved(int)'
$ws11.Cells.Item(22,3).Value = '3'

$ws11.Cells.Item(23,1).Value = 'pl.piomin.payment.domain.Customer'
$ws11.Cells.Item(23,2).Value = 'toString()'
$ws11.Cells.Item(23,3).Value = '3'

$ws11.Cells.Item(24,1).Value = 'pl.piomin.payment.KafkaContainerDevMode'
$ws11.Cells.Item(24,2).Value = 'KafkaContainerDevMode()'
$ws11.Cells.Item(24,3).Value = '1'

$ws11.Cells.Item(25,1).Value = 'pl.piomin.payment.KafkaContainerDevMode'
$ws11.Cells.Item(25,2).Value = 'kafka()'
$ws11.Cells.Item(25,3).Value = '3'

